# Fix Training Data Issue: the "Date" column was populated with the
# source filename/label ("6-22-2007-08") instead of the actual game
# date. Correct it to the real ISO date (2008-06-22) for every data row,
# while keeping the cell as literal text (not re-interpreted as an
# Excel date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Locate the "Date" header column on row 1.
$dateCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value() -eq "Date") {
        $dateCol = $c
    }
}

$oldValue = "6-22-2007-08"
$newValue = "2008-06-22"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    if ($cell.Value() -eq $oldValue) {
        # Assigning the literal string via .Value would make Excel's
        # type-inference re-parse "2008-06-22" as a date serial number.
        # Building it as a text formula and collapsing it back to a
        # static value with Paste Special (values only) keeps it a
        # plain text cell, matching the original cell's data type/style.
        $cell.Formula = "=""" + $newValue + """"
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = 0
